$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.998.22"
$ws.Range("E2").Value = "  -0.52%  "

$ws.Range("D3").Value = "2.037.57"
$ws.Range("E3").Value = "  -1.02%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.74"
$ws.Range("E5").Value = "  -0.84%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.660"
$ws.Range("E6").Value = "  -0.58%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.95"
$ws.Range("E8").Value = "  +1.59%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.378"
$ws.Range("E9").Value = "  -1.36%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0774"
$ws.Range("E10").Value = "  +0.98%  "

$ws.Range("E11").Value = "  +0.96%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.75"
$ws.Range("E12").Value = "  +3.19%  "

$ws.Range("D13").Value = "2.331.66"
$ws.Range("E13").Value = "  -1.08%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.58"
$ws.Range("E14").Value = "  +4.88%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.789"
$ws.Range("E15").Value = "  -5.63%  "

$ws.Range("D16").Value = "2.031.30"
$ws.Range("E16").Value = "  -1.31%  "

$ws.Range("D17").Value = "36.916.13"
$ws.Range("E17").Value = "  -0.27%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "16.47"
$ws.Range("E18").Value = "  +12.89%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.93"
$ws.Range("E19").Value = "  +1.32%  "

$ws.Range("D20").Value = "0.0₃0889"
$ws.Range("E20").Value = "  +3.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.29"
$ws.Range("E21").Value = "  +0.00%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "234.98"
$ws.Range("E22").Value = "  -1.78%  "

$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.36"
$ws.Range("E24").Value = "  -3.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.24"
$ws.Range("E25").Value = "  +10.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.09"
$ws.Range("E26").Value = "  -1.84%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.09"
$ws.Range("E27").Value = "  -1.25%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.70"
$ws.Range("E28").Value = "  -4.82%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.124"
$ws.Range("E29").Value = "  +0.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.12"
$ws.Range("E30").Value = "  +0.89%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.66"
$ws.Range("E31").Value = "  +1.26%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0609"
$ws.Range("E32").Value = "  -3.98%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.38"
$ws.Range("E33").Value = "  -0.09%  "

$ws.Range("E34").Value = "  -0.06%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0865"
$ws.Range("E35").Value = "  -4.08%  "

$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.19"
$ws.Range("E36").Value = "  -4.35%  "

$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.76"
$ws.Range("E37").Value = "  -3.79%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.33"
$ws.Range("E38").Value = "  -1.67%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.105"
$ws.Range("E39").Value = "  +0.40%  "

$ws.Range("B40").Value = "HuobiToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.07"
$ws.Range("E40").Value = "  +9.68%  "

$ws.Range("B41").Value = "THORChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.90"
$ws.Range("E41").Value = "  +22.14%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0219"
$ws.Range("E42").Value = "  -3.81%  "

$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.27"
$ws.Range("E43").Value = "  -5.35%  "

$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "94.88"
$ws.Range("E44").Value = "  -2.93%  "

$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.10"
$ws.Range("E45").Value = "  -4.96%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.41"
$ws.Range("E46").Value = "  -0.66%  "

$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D47").Value = "1.271.11"
$ws.Range("E47").Value = "  -3.35%  "

$ws.Range("B48").Value = "MXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.83"
$ws.Range("E48").Value = "  -2.82%  "

$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.66"
$ws.Range("E49").Value = "  -4.82%  "

$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.220.53"
$ws.Range("E50").Value = "  -0.76%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "42.28"
$ws.Range("E51").Value = "  -7.11%  "
